# menambahkan logika sheet billing
#
# - insert a new "Weekly" sheet (gets the data currently on "Main")
# - clear "Main" so it becomes empty
# - append a new "Billing" sheet (after "Monthly") with fresh data
#
$wb = $excel.ActiveWorkbook

# --- 1. "Weekly" sheet: carries the data that used to live on "Main" ---
$main  = $wb.Worksheets.Item("Main")
$daily = $wb.Worksheets.Item("Daily")

$weekly = $wb.Worksheets.Add($null, $daily)
$weekly.Name = "Weekly"

$main.Range("A1:E5").Copy()
$weekly.Range("A1").PasteSpecial()

# --- 2. empty out "Main" now that its data moved to "Weekly" ---
$main.Cells.Clear()

# --- 3. "Billing" sheet: brand-new data, placed after "Monthly" ---
# (re-fetch "Monthly" now, since inserting "Weekly" shifted its position and
#  a handle obtained earlier would now resolve to the wrong sheet)
$monthly = $wb.Worksheets.Item("Monthly")
$billing = $wb.Worksheets.Add($null, $monthly)
$billing.Name = "Billing"

$billing.Range("A1").Value = "ih_biling_ih"

$billing.Range("A2").Value = "TABLE NAME"
$billing.Range("B2").Value = "EVENT DATE"
$billing.Range("C2").Value = "DATE TRANSACTION"
$billing.Range("D2").Value = "DATE AVAILABILITY"
$billing.Range("E2").Value = "NOW SIZE CONDITION"

$billing.Range("A3").Value = "ih_biling_ih"
$billing.Range("B3").Value = "event_date=2024-11-05"
# force as text: a bare "2024-11-06" would otherwise be auto-parsed as a date
$billing.Range("C3").NumberFormat = "@"
$billing.Range("C3").Value = "2024-11-06"
$billing.Range("D3").Value = "04:50"
$billing.Range("E3").Value = "36.6 K"
